$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 321, shifting existing rows 321:391 down to 322:392
$ws.Rows("321:321").Insert()

# Populate the newly inserted row 321 with a new weekly entry.
# (Same attributes as the former row 321, but with an updated, more recent date.)
$ws.Cells.Item(321, 1).Value = 5
$ws.Cells.Item(321, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(321, 3).Value = "Maule"
$ws.Range("D321").Value = 44995
$ws.Cells.Item(321, 5).Value = 7
$ws.Cells.Item(321, 6).Value = 100112009
$ws.Cells.Item(321, 7).Value = "Acelga"
$ws.Cells.Item(321, 8).Value = "Sin especificar"
$ws.Cells.Item(321, 9).Value = "Primera"
$ws.Cells.Item(321, 10).Value = 500
$ws.Cells.Item(321, 11).Value = 2500
$ws.Cells.Item(321, 12).Value = 2500
$ws.Cells.Item(321, 13).Value = 2500
$ws.Cells.Item(321, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(321, 15).Value = "Región del Maule"
$ws.Cells.Item(321, 16).Value = 625
$ws.Cells.Item(321, 17).Value = 4
$ws.Cells.Item(321, 18).Value = "Hortaliza"
